# Update crypto tracker values per latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.207.71'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.890.49'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '242.57'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '0.653'
$ws.Range('E6').Value = '  +5.15%  '
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('D8').Value = '41.09'
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').Value = '0.346'
$ws.Range('E9').Value = '  +5.87%  '
$ws.Range('D10').Value = '50.07'
$ws.Range('E10').Value = '  +7.63%  '
$ws.Range('D11').Value = '0.0706'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '2.167.27'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').Value = '11.87'
$ws.Range('E14').Value = '  +4.30%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.897.91'
$ws.Range('E15').Value = '  +3.31%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.690'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('D18').Value = '35.213.28'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '71.05'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = '0.0₃0809'
$ws.Range('E20').Value = '  +2.13%  '
$ws.Range('D21').Value = '240.31'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '12.36'
$ws.Range('E22').Value = '  +1.44%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('D25').Value = '2.43'
$ws.Range('E25').Value = '  +32.82%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('D27').Value = '169.91'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('D28').Value = '8.39'
$ws.Range('E28').Value = '  +5.21%  '
$ws.Range('D29').Value = '18.16'
$ws.Range('E29').Value = '  +3.26%  '
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('D31').Value = '4.09'
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('B32').Value = 'BinanceUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D32').Value = '1.01'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0558'
$ws.Range('E33').Value = '  +0.67%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.935'
$ws.Range('E34').Value = '  +15.84%  '
$ws.Range('D35').Value = '4.08'
$ws.Range('E35').Value = '  +1.85%  '
$ws.Range('E36').Value = '  -1.69%  '
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('D38').Value = '1.31'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('E39').Value = '  +3.54%  '
$ws.Range('E40').Value = '  +1.43%  '
$ws.Range('E41').Value = '  +15.15%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.90'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +7.09%  '
$ws.Range('D43').Value = '88.75'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('D44').Value = '1.335.02'
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = '2.35'
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('B46').Value = 'MultiversX'
$ws.Range('C46').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D46').Value = '47.81'
$ws.Range('E46').Value = '  +38.62%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.40'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('D48').Value = '2.76'
$ws.Range('E48').Value = '  +0.77%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').Value = '2.076.83'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').Value = '11.16'
$ws.Range('E51').Value = '  -14.00%  '
